# Applies the 2021-01-05 doc edit: trims two intro paragraphs, splits the
# date line into "Date :" / " 2021 - 01- 05", and sprinkles <w:proofErr/>
# spell/grammar-check markers (plus a couple of page-break fixes) across
# several paragraphs, matching Word's own re-save behaviour.

$d = $word.ActiveDocument

function Get-Xml([string]$innerRunsXml) {
    return '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<w:body><w:p>' + $innerRunsXml + '</w:p></w:body></w:document>' +
        '</pkg:xmlData></pkg:part></pkg:package>'
}

# Replace the exact text $needle (first match, plain text) with the raw
# run/proofErr markup $innerRunsXml, preserving the paragraph it lives in.
function Replace-RunText([string]$needle, [string]$innerRunsXml) {
    $rng = $d.Content
    $found = $rng.Find.Execute($needle, $true, $true, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        throw "Find failed for: $needle"
    }
    $target = $d.Range($rng.Start, $rng.End)
    $target.InsertXML((Get-Xml $innerRunsXml))
}

$rPr = '<w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-US"/></w:rPr>'

# --- 1) Delete "This is first change " and "This is second" paragraphs ---
$rng = $d.Content
$found = $rng.Find.Execute("This is first change ", $true, $true, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) { throw "Find failed for first-change paragraph" }
$d.Range($rng.Start, $rng.End + 1).Delete()

$rng = $d.Content
$found = $rng.Find.Execute("This is second", $true, $true, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) { throw "Find failed for second paragraph" }
$d.Range($rng.Start, $rng.End + 1).Delete()

# --- 2) "Date : 2021 - 01- 05" -> "Date :" (gram-checked) + " 2021 - 01- 05" ---
$inner = '<w:proofErr w:type="gramStart"/>' +
         '<w:r>' + $rPr + '<w:t>Date :</w:t></w:r>' +
         '<w:proofErr w:type="gramEnd"/>' +
         '<w:r>' + $rPr + '<w:t xml:space="preserve"> 2021 - 01- 05</w:t></w:r>'
Replace-RunText "Date : 2021 - 01- 05" $inner

# --- 3) "Create a account on github." -> split w/ proofErr around "a" and "github" ---
$inner = '<w:r>' + $rPr + '<w:t xml:space="preserve">Create </w:t></w:r>' +
         '<w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/>' +
         '<w:r>' + $rPr + '<w:t>a</w:t></w:r>' +
         '<w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/>' +
         '<w:r>' + $rPr + '<w:t xml:space="preserve"> account on </w:t></w:r>' +
         '<w:proofErr w:type="spellStart"/>' +
         '<w:r>' + $rPr + '<w:t>github</w:t></w:r>' +
         '<w:proofErr w:type="spellEnd"/>' +
         '<w:r>' + $rPr + '<w:t>.</w:t></w:r>'
Replace-RunText "Create a account on github." $inner

# --- 4) ssh-keygen line -> split w/ proofErr around ssh / rsa / piparava ---
$quote = [char]0x201c
$rquote = [char]0x201d
$needle = "ssh-keygen -t rsa -C " + $quote + "jay-piparava" + $rquote
$inner = '<w:proofErr w:type="spellStart"/>' +
         '<w:r>' + $rPr + '<w:t>ssh</w:t></w:r>' +
         '<w:proofErr w:type="spellEnd"/>' +
         '<w:r>' + $rPr + '<w:t xml:space="preserve">-keygen -t </w:t></w:r>' +
         '<w:proofErr w:type="spellStart"/>' +
         '<w:r>' + $rPr + '<w:t>rsa</w:t></w:r>' +
         '<w:proofErr w:type="spellEnd"/>' +
         '<w:r>' + $rPr + '<w:t xml:space="preserve"> -C ' + $quote + 'jay-</w:t></w:r>' +
         '<w:proofErr w:type="spellStart"/>' +
         '<w:r>' + $rPr + '<w:t>piparava</w:t></w:r>' +
         '<w:proofErr w:type="spellEnd"/>' +
         '<w:r>' + $rPr + '<w:t>' + $rquote + '</w:t></w:r>'
Replace-RunText $needle $inner

# --- 5) "copy ssh key and add it." -> split w/ proofErr around "ssh" ---
$inner = '<w:r>' + $rPr + '<w:t xml:space="preserve">copy </w:t></w:r>' +
         '<w:proofErr w:type="spellStart"/>' +
         '<w:r>' + $rPr + '<w:t>ssh</w:t></w:r>' +
         '<w:proofErr w:type="spellEnd"/>' +
         '<w:r>' + $rPr + '<w:t xml:space="preserve"> key and add it.</w:t></w:r>'
Replace-RunText "copy ssh key and add it." $inner

# --- 6) "mkdir [name] " -> split w/ proofErr around "mkdir" ---
$inner = '<w:proofErr w:type="spellStart"/>' +
         '<w:r>' + $rPr + '<w:t>mkdir</w:t></w:r>' +
         '<w:proofErr w:type="spellEnd"/>' +
         '<w:r>' + $rPr + '<w:t xml:space="preserve"> [name] </w:t></w:r>'
Replace-RunText "mkdir [name] " $inner

# --- 7) "now create a branch so that a error ... can not affected with error." ---
$inner = '<w:r>' + $rPr + '<w:t xml:space="preserve">now create a branch so that </w:t></w:r>' +
         '<w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/>' +
         '<w:r>' + $rPr + '<w:t>a</w:t></w:r>' +
         '<w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/>' +
         '<w:r>' + $rPr + '<w:t xml:space="preserve"> error are detected in branch so main repository </w:t></w:r>' +
         '<w:proofErr w:type="spellStart"/>' +
         '<w:r>' + $rPr + '<w:t>can not</w:t></w:r>' +
         '<w:proofErr w:type="spellEnd"/>' +
         '<w:r>' + $rPr + '<w:t xml:space="preserve"> affected with error.</w:t></w:r>'
Replace-RunText "now create a branch so that a error are detected in branch so main repository can not affected with error." $inner

# --- 8) "for adding a file use commad" -> split w/ proofErr around "commad" ---
$inner = '<w:r>' + $rPr + '<w:t xml:space="preserve">for adding a file use </w:t></w:r>' +
         '<w:proofErr w:type="spellStart"/>' +
         '<w:r>' + $rPr + '<w:t>commad</w:t></w:r>' +
         '<w:proofErr w:type="spellEnd"/>'
Replace-RunText "for adding a file use commad" $inner

# --- 9) Move <w:lastRenderedPageBreak/> from "now commit changes " to "now push changes..." ---
$inner = '<w:r>' + $rPr + '<w:t xml:space="preserve">now commit changes </w:t></w:r>'
Replace-RunText "now commit changes " $inner

$inner = '<w:r>' + $rPr + '<w:lastRenderedPageBreak/><w:t>now push changes in original repository</w:t></w:r>'
Replace-RunText "now push changes in original repository" $inner

# --- 10) "git push origin[branch name]" -> split w/ proofErr around "origin[" ---
$inner = '<w:r>' + $rPr + '<w:t xml:space="preserve">git push </w:t></w:r>' +
         '<w:proofErr w:type="gramStart"/>' +
         '<w:r>' + $rPr + '<w:t>origin[</w:t></w:r>' +
         '<w:proofErr w:type="gramEnd"/>' +
         '<w:r>' + $rPr + '<w:t>branch name]</w:t></w:r>'
Replace-RunText "git push origin[branch name]" $inner

Write-Output "done"
